$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: replace organism name with quoted .mat filename literal.
# A leading "'" is Excel's text-qualifier prefix and gets stripped from the
# stored value, so we double it up ("''...") to end up with a single literal
# leading quote in the saved text. Re-applying the "Normal" style afterwards
# clears the auto-added quote-prefix cell style so formatting stays untouched.
# Column C: updated numeric value (rounded fractions from the new source data).

$cell = $ws.Cells.Item(2, 2)
$cell.Value = "''Akkermansia_muciniphila_ATCC_BAA_835.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(2, 3).Value = 0

$cell = $ws.Cells.Item(3, 2)
$cell.Value = "''Alistipes_finegoldii_DSM_17242.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(3, 3).Value = 0

$cell = $ws.Cells.Item(4, 2)
$cell.Value = "''Alistipes_indistinctus_YIT_12060.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(4, 3).Value = 0

$cell = $ws.Cells.Item(5, 2)
$cell.Value = "''Alistipes_putredinis_DSM_17216.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(5, 3).Value = 0

$cell = $ws.Cells.Item(6, 2)
$cell.Value = "''Alistipes_shahii_WAL_8301.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(6, 3).Value = 0

$cell = $ws.Cells.Item(7, 2)
$cell.Value = "''Bacteroides_cellulosilyticus_DSM_14838.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(7, 3).Value = 0.004

$cell = $ws.Cells.Item(8, 2)
$cell.Value = "''Bacteroides_coprophilus_DSM_18228.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(8, 3).Value = 0.001

$cell = $ws.Cells.Item(9, 2)
$cell.Value = "''Bacteroides_fragilis_3_1_12.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(9, 3).Value = 0.017

$cell = $ws.Cells.Item(10, 2)
$cell.Value = "''Bacteroides_oleiciplenus_YIT_12058.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(10, 3).Value = 0

$cell = $ws.Cells.Item(11, 2)
$cell.Value = "''Bacteroides_ovatus_ATCC_8483.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(11, 3).Value = 0.089

$cell = $ws.Cells.Item(12, 2)
$cell.Value = "''Bacteroides_plebeius_M12_DSM_17135.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(12, 3).Value = 0

$cell = $ws.Cells.Item(13, 2)
$cell.Value = "''Bacteroides_salyersiae_WAL_10018.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(13, 3).Value = 0.047

$cell = $ws.Cells.Item(14, 2)
$cell.Value = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(14, 3).Value = 0

$cell = $ws.Cells.Item(15, 2)
$cell.Value = "''Bacteroides_uniformis_ATCC_8492.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(15, 3).Value = 0

$cell = $ws.Cells.Item(16, 2)
$cell.Value = "''Bacteroides_vulgatus_ATCC_8482.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(16, 3).Value = 0.049

$cell = $ws.Cells.Item(17, 2)
$cell.Value = "''Barnesiella_intestinihominis_YIT_11860.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(17, 3).Value = 0

$cell = $ws.Cells.Item(18, 2)
$cell.Value = "''Bifidobacterium_animalis_lactis_AD011.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(18, 3).Value = 0

$cell = $ws.Cells.Item(19, 2)
$cell.Value = "''Bilophila_wadsworthia_3_1_6.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(19, 3).Value = 0

$cell = $ws.Cells.Item(20, 2)
$cell.Value = "''Escherichia_coli_O157_H7_str_Sakai.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(20, 3).Value = 0.045

$cell = $ws.Cells.Item(21, 2)
$cell.Value = "''Eubacterium_limosum_KIST612.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(21, 3).Value = 0.748

$cell = $ws.Cells.Item(22, 2)
$cell.Value = "''Eubacterium_ramulus_ATCC_29099.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(22, 3).Value = 0

$cell = $ws.Cells.Item(23, 2)
$cell.Value = "''Flavonifractor_plautii_ATCC_29863.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(23, 3).Value = 0

$cell = $ws.Cells.Item(24, 2)
$cell.Value = "''Marvinbryantia_formatexigens_I_52_DSM_14469.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(24, 3).Value = 0

$cell = $ws.Cells.Item(25, 2)
$cell.Value = "''Odoribacter_splanchnicus_1651_6_DSM_20712.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(25, 3).Value = 0

$cell = $ws.Cells.Item(26, 2)
$cell.Value = "''Parabacteroides_distasonis_ATCC_8503.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(26, 3).Value = 0

$cell = $ws.Cells.Item(27, 2)
$cell.Value = "''Parabacteroides_johnsonii_DSM_18315.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(27, 3).Value = 0

$cell = $ws.Cells.Item(28, 2)
$cell.Value = "''Paraprevotella_xylaniphila_YIT_11841.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(28, 3).Value = 0

$cell = $ws.Cells.Item(29, 2)
$cell.Value = "''Parasutterella_excrementihominis_YIT_11859.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(29, 3).Value = 0

$cell = $ws.Cells.Item(30, 2)
$cell.Value = "''Phascolarctobacterium_succinatutens_YIT_12067.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(30, 3).Value = 0

$cell = $ws.Cells.Item(31, 2)
$cell.Value = "''Prevotella_copri_CB7_DSM_18205.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(31, 3).Value = 0

$cell = $ws.Cells.Item(32, 2)
$cell.Value = "''Prevotella_stercorea_DSM_18206.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(32, 3).Value = 0

$cell = $ws.Cells.Item(33, 2)
$cell.Value = "''Roseburia_inulinivorans_DSM_16841.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(33, 3).Value = 0

$cell = $ws.Cells.Item(34, 2)
$cell.Value = "''Sutterella_wadsworthensis_3_1_45B.mat'"
$cell.Style = "Normal"
$ws.Cells.Item(34, 3).Value = 0
